# Updates cryptos list data (prices & 1h volume %) per the GitHub Actions commit.
# Rows 12 and 13 additionally swap their Coin/Link/Price/Volume content (ranking reorder).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.804.46'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.902.08'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.94'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5017'
$ws.Range("E7").Value = '  +3.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3807'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07275'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9078'
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07651'
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.922.14'
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.481'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.24'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008693'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.842.95'
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.163'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.22'
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.857'
$ws.Range("E25").Value = '  -3.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.228'
$ws.Range("E26").Value = '  +5.84%  '
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.22'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.905'
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08963'
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.205'
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.232'
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7658'
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.637'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02057'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.545'
$ws.Range("E36").Value = '  -2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.096'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5546'
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.013'
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05252'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.970'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.482'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1522'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '110.86'
$ws.Range("E44").Value = '  +3.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.61'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4785'
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.630'
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.25'
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06073'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8991'
$ws.Range("E51").Value = '  -0.37%  '
